$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Colour the word "Subject" (first occurrence, in the Subject line)
#    red. Word will split the run so only "Subject" carries the new
#    <w:color w:val="FF0000"/>, the remaining text keeps its original
#    (colourless) run.
# ------------------------------------------------------------------
$rngSubject = $d.Content
$rngSubject.Find.Execute("Subject", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngSubject.Font.Color = 255

# ------------------------------------------------------------------
# 2) Colour the word "Requirements" red as well.
# ------------------------------------------------------------------
$rngReq = $d.Content
$rngReq.Find.Execute("Requirements", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngReq.Font.Color = 255

# ------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the "User guide" /
#    "recipes" paragraph to the end of the paragraph that ends with
#    "... add or remove navigation properties."
# ------------------------------------------------------------------

# Remove the bookmark from its old location (right after "recipes").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the end of "... add or remove navigation properties."
$rngNav = $d.Content
$rngNav.Find.Execute("add or remove navigation properties.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngNav.Collapse(0)

# A zero-length range collapsed exactly on a paragraph mark cannot be
# handed straight to Bookmarks.Add and land in the right paragraph, so
# insert a tiny placeholder, bookmark around it, then delete the
# placeholder text again - this leaves bookmarkStart/bookmarkEnd as the
# last two children of the paragraph, exactly like the original markup.
$rngNav.InsertAfter("\u0001GoBackMark\u0001")

$rngMark = $d.Content
$rngMark.Find.Execute("\u0001GoBackMark\u0001", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $rngMark)

$rngMark2 = $d.Content
$rngMark2.Find.Execute("\u0001GoBackMark\u0001", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngMark2.Text = ""
